$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7090943455696106
$ws.Range("B1").Value = 1.975718021392822
$ws.Range("C1").Value = 2.503110408782959
$ws.Range("D1").Value = 0.8493666052818298
$ws.Range("E1").Value = 0.958307147026062
